# DeliveryChanges.xlsx update
#
# Changes applied (per the target diff):
#  1. Row 13: remove the Start Date value in A13 (cell becomes blank,
#     not just value-cleared but removed entirely along with its format).
#  2. Row 16: remove the Detail value in C16 ("Frozen in nose from 4-07
#     to 4-28") - cell removed entirely.
#  3. Remove rows 22 and 23 completely (SA 16 / Inventory on 04-15 row,
#     and the stray "123/4" row), shifting nothing else since they were
#     the last two rows.
#
# Using .Clear() (not ClearContents) so the cell's style is dropped too -
# this makes the cell fully empty/absent, matching the target XML where
# the <c> element disappears rather than being left as an empty tag.
# Using Rows(...).Delete() removes the rows outright (shifting the used
# range / dimension up), which also lets now-unreferenced shared strings
# get garbage-collected on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the Start Date cell in row 13 (A13) entirely.
$ws.Range("A13").Clear()

# 2. Clear the Detail cell in row 16 (C16) entirely.
$ws.Range("C16").Clear()

# 3. Delete rows 22 and 23 entirely (last two rows of the sheet).
$ws.Rows("22:23").Delete()
